$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two extra biopsies sequenced on a new chip ("1234567"), reported later
# (25-may-2023). Row 10 mirrors the "clean" pattern (row 2); row 11 mirrors
# the KRAS-positive pattern (row 3).
#
# Columns A, B, C, E hold numeric-looking codes and F holds a date-looking
# string, but — like the rest of the sheet — they must stay plain text, so
# format those cells as Text before writing the values (keeps Excel from
# auto-converting "100"/"1"/"1234567"/"25-may-2023" into a number/date).
$textCols = @("A10:A11", "B10:B11", "C10:C11", "E10:E11", "F10:F11")
foreach ($addr in $textCols) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 10
$ws.Cells.Item(10, 1).Value = "100"
$ws.Cells.Item(10, 2).Value = "1"
$ws.Cells.Item(10, 3).Value = "1234567"
$ws.Cells.Item(10, 4).Value = "23B000000-A1/CHIP100.1"
$ws.Cells.Item(10, 5).Value = "1"
$ws.Cells.Item(10, 6).Value = "25-may-2023"
$ws.Cells.Item(10, 7).Value = "Carcinoma pulmonar no microcítico"
$ws.Cells.Item(10, 8).Value = 15.1
$ws.Cells.Item(10, 9).Value = "[]"
$ws.Cells.Item(10, 10).Value = "[]"
$ws.Cells.Item(10, 11).Value = "[]"
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 6
$ws.Cells.Item(10, 14).Value = 1
$ws.Cells.Item(10, 15).Value = 6
$ws.Cells.Item(10, 16).Value = 1

# Row 11
$ws.Cells.Item(11, 1).Value = "100"
$ws.Cells.Item(11, 2).Value = "2"
$ws.Cells.Item(11, 3).Value = "1234567"
$ws.Cells.Item(11, 4).Value = "23B00000-A1/CHIP100.2"
$ws.Cells.Item(11, 5).Value = "1"
$ws.Cells.Item(11, 6).Value = "25-may-2023"
$ws.Cells.Item(11, 7).Value = "Carcinoma pulmonar no microcítico"
$ws.Cells.Item(11, 8).Value = 15.1
$ws.Cells.Item(11, 9).Value = "['KRAS']"
$ws.Cells.Item(11, 10).Value = "[35]"
$ws.Cells.Item(11, 11).Value = "['66.50']"
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 4
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = 1
$ws.Cells.Item(11, 16).Value = 1
